$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (245-247), continuing the existing series in A:D.
$data = @(
    @(44319, 4, 74, 183.7231242862108),
    @(44320, 5, 70, 173.7921445950643),
    @(44321, 6, 73, 181.2403793634242)
)

$startRow = 245
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Column A carries the date style (s="2") used by the rest of the
    # column; replicate it from the row above via a format-only paste
    # so we reuse the existing style entry instead of fabricating a new
    # (slightly different) one through individual font/format writes.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
